$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Model_name" in E1, matching the style of the existing
# header cells (A1:D1) by copying their formatting.
$ws.Range("E1").Value = "Model_name"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the re-fit peak values in row 2
$ws.Range("B2").Value = 1089.21919294995
$ws.Range("C2").Value = 763.4596414328784
$ws.Range("D2").Value = 50.18331271259285

# Add the new model-name value for row 2
$ws.Range("E2").Value = "Spline"
